# Add one more UWS sample row (row 7) to the continuous file list, mirroring
# the formatting already used by the previous UWS rows (rows 4-6), and
# re-center the "flag" column (E) and vertically-center the "#" column (B)
# for all of the data rows, including the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7 values ---------------------------------------------------
$ws.Range("A7").Value = "UWS"
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "2022-02-24_221145_SO289"
$ws.Range("E7").Value = 2

# --- Formatting -----------------------------------------------------------
# "flag" column (E): was left-aligned, now centered for every data row,
# including the newly added one.
$ws.Range("E3:E7").HorizontalAlignment = -4108

# "sample #" column (B): center horizontally and vertically for every data
# row, including the newly added one.
$ws.Range("B3:B7").HorizontalAlignment = -4108
$ws.Range("B3:B7").VerticalAlignment = -4108

# "sample" column (A) on the new row matches the centered style used by the
# other UWS rows above it.
$ws.Range("A7").HorizontalAlignment = -4108

# --- Update the selection to reflect the newly edited range --------------
$null = $ws.Range("E3:E7").Select()
